$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.324.11"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.091.05"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.083.78"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  +5.37%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.74%  "
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "3.600.61"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "63.184.75"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "3.087.30"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.82%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").Value = "0.0₃0843"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "434.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0366"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "2.863.56"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.70%  "
